$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.030.63"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.508.86"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'589.30"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "'177.07"
$ws.Range("E6").Value = "  +3.66%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.74%  "
$ws.Range("E9").Value = "  +5.08%  "
$ws.Range("E10").Value = "  -0.71%  "
$ws.Range("D11").Value = "'0.340"
$ws.Range("E11").Value = "  +1.91%  "
$ws.Range("D12").Value = "'4.96"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").Value = "2.978.99"
$ws.Range("E13").Value = "  +3.16%  "
$ws.Range("D14").Value = "'25.76"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "67.842.29"
$ws.Range("E15").Value = "  +1.19%  "
$ws.Range("D16").Value = "'0.0000172"
$ws.Range("E16").Value = "  +1.13%  "
$ws.Range("D17").Value = "2.510.45"
$ws.Range("E17").Value = "  +5.04%  "
$ws.Range("D18").Value = "'11.00"
$ws.Range("E18").Value = "  -0.24%  "
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").Value = "'353.03"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("D21").Value = "'4.12"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "'70.99"
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("D24").Value = "'4.31"
$ws.Range("E24").Value = "  +1.75%  "
$ws.Range("D25").Value = "'1.75"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "'9.20"
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "2.638.06"
$ws.Range("E27").Value = "  +0.60%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'509.57"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("E31").Value = "  +0.99%  "
$ws.Range("E32").Value = "  +2.68%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "'164.61"
$ws.Range("E35").Value = "  +2.60%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.121"
$ws.Range("E36").Value = "  +2.54%  "
$ws.Range("D37").Value = "'18.41"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("D39").Value = "'1.34"
$ws.Range("E39").Value = "  +0.10%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("E42").Value = "  +0.38%  "
$ws.Range("D43").Value = "'4.88"
$ws.Range("E43").Value = "  +1.68%  "
$ws.Range("D44").Value = "'2.50"
$ws.Range("E44").Value = "  +5.76%  "
$ws.Range("D45").Value = "'147.72"
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("D46").Value = "'3.57"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("D47").Value = "0.0₆0262"
$ws.Range("E47").Value = "  +4.23%  "
$ws.Range("E48").Value = "  +1.17%  "
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("D51").Value = "'0.587"
$ws.Range("E51").Value = "  +0.65%  "

Write-Output "done"
